$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.416.11"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "'3.563.04"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'587.60"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'187.78"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("D7").Value = "'3.556.35"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.216"
$ws.Range("E10").Value = "  +9.52%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "'0.0000312"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "'9.42"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "'4.128.84"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'3.606.02"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'70.433.02"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "'12.76"
$ws.Range("E18").Value = "  +3.59%  "
$ws.Range("D19").Value = "'19.00"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").Value = "'578.67"
$ws.Range("E20").Value = "  +8.74%  "
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "'17.67"
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("D24").Value = "'4.64"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "'4.89"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "'94.45"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'2.93"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.96"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").Value = "'9.37"
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("D30").Value = "'32.31"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").Value = "'7.08"
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("D32").Value = "'12.23"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "'0.114"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").Value = "'63.86"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'3.77"
$ws.Range("E35").Value = "  +22.17%  "
$ws.Range("D36").Value = "'3.22"
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("D37").Value = "'3.831.45"
$ws.Range("E37").Value = "  +14.39%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "'38.12"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'521.52"
$ws.Range("E40").Value = "  -4.61%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'0.0₃0788"
$ws.Range("E42").Value = "  +3.85%  "
$ws.Range("D43").Value = "'3.57"
$ws.Range("E43").Value = "  +5.65%  "
$ws.Range("E44").Value = "  +2.59%  "
$ws.Range("D45").Value = "'0.0454"
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("D46").Value = "'2.95"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "'3.46"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("E51").Value = "  +6.80%  "
